$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for columns D and E (prices / percentages are
# stored as literal strings, e.g. "29.058.02", "0.9990") so Excel does not
# silently reinterpret them as numbers and lose the exact text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.058.02"

$ws.Range("D3").Value = "1.828.72"

$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "241.37"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "0.6366"
$ws.Range("E6").Value = "  -4.22%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "44.73"
$ws.Range("E8").Value = "  +6.66%  "

$ws.Range("D9").Value = "0.2937"
$ws.Range("E9").Value = "  +0.78%  "

$ws.Range("D10").Value = "0.07333"
$ws.Range("E10").Value = "  -0.28%  "

$ws.Range("E11").Value = "  +0.92%  "

$ws.Range("D12").Value = "0.07653"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").Value = "1.827.27"
$ws.Range("E13").Value = "  -0.36%  "

$ws.Range("D14").Value = "4.986"

$ws.Range("D15").Value = "0.6640"
$ws.Range("E15").Value = "  -0.05%  "

$ws.Range("D16").Value = "81.99"
$ws.Range("E16").Value = "  -1.73%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "6.054"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000008655"
$ws.Range("E18").Value = "  +4.87%  "

$ws.Range("D19").Value = "28.904.03"
$ws.Range("E19").Value = "  -0.99%  "

$ws.Range("D20").Value = "2.077.32"
$ws.Range("E20").Value = "  -0.31%  "

$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").Value = "224.11"
$ws.Range("E22").Value = "  -0.61%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").Value = "7.116"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").Value = "1.000"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("D26").Value = "157.78"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("D27").Value = "8.464"
$ws.Range("E27").Value = "  -1.84%  "

$ws.Range("D28").Value = "0.1372"
$ws.Range("E28").Value = "  -1.31%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").Value = "1.502"
$ws.Range("E30").Value = "  -0.27%  "

$ws.Range("D31").Value = "4.095"
$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("D32").Value = "4.027"
$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E33").Value = "  +1.55%  "

$ws.Range("D34").Value = "0.05290"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("D35").Value = "1.836"
$ws.Range("E35").Value = "  -1.52%  "

$ws.Range("D36").Value = "0.7378"
$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("D37").Value = "1.153"
$ws.Range("E37").Value = "  +2.16%  "

$ws.Range("D38").Value = "2.653"
$ws.Range("E38").Value = "  -0.92%  "

$ws.Range("D39").Value = "1.292.72"
$ws.Range("E39").Value = "  -0.69%  "

$ws.Range("D40").Value = "2.752"
$ws.Range("E40").Value = "  +1.28%  "

$ws.Range("D41").Value = "0.01782"
$ws.Range("E41").Value = "  -0.76%  "

$ws.Range("D42").Value = "6.293"
$ws.Range("E42").Value = "  +6.06%  "

$ws.Range("D43").Value = "0.8945"
$ws.Range("E43").Value = "  -2.70%  "

$ws.Range("E44").Value = "  -0.13%  "

$ws.Range("D45").Value = "102.90"
$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").Value = "1.975.96"
$ws.Range("E46").Value = "  -0.28%  "

$ws.Range("D47").Value = "0.5135"
$ws.Range("E47").Value = "  -0.51%  "

$ws.Range("D48").Value = "63.93"
$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("E49").Value = "  -6.16%  "

$ws.Range("D50").Value = "1.730"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").Value = "0.07237"
$ws.Range("E51").Value = "  -14.90%  "

